$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.901.58"
$ws.Range("E2").Value = "  +3.95%  "

$ws.Range("D3").Value = "2.652.19"
$ws.Range("E3").Value = "  +6.15%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.71"
$ws.Range("E5").Value = "  +2.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.07"
$ws.Range("E6").Value = "  +3.49%  "

$ws.Range("E7").Value = "  +1.34%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.559"
$ws.Range("E9").Value = "  +3.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.93"
$ws.Range("E10").Value = "  +3.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.48"
$ws.Range("E11").Value = "  +1.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0822"
$ws.Range("E12").Value = "  +1.46%  "

$ws.Range("E13").Value = "  +0.77%  "

$ws.Range("E14").Value = "  +3.47%  "

$ws.Range("D15").Value = "3.074.31"
$ws.Range("E15").Value = "  +6.37%  "

$ws.Range("D16").Value = "2.687.24"
$ws.Range("E16").Value = "  +7.38%  "

$ws.Range("E17").Value = "  +5.90%  "

$ws.Range("D18").Value = "49.889.04"
$ws.Range("E18").Value = "  +4.20%  "

$ws.Range("E19").Value = "  +3.19%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.83"
$ws.Range("E20").Value = "  +1.97%  "

$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.95"
$ws.Range("E21").Value = "  +7.22%  "

$ws.Range("D22").Value = "0.0₃0960"
$ws.Range("E22").Value = "  +2.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.14"
$ws.Range("E23").Value = "  +2.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "280.79"
$ws.Range("E24").Value = "  +1.81%  "

$ws.Range("E25").Value = "  +2.82%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.06"
$ws.Range("E26").Value = "  +4.79%  "

$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.24"
$ws.Range("E28").Value = "  +2.03%  "

$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.88"
$ws.Range("E29").Value = "  +4.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.98"
$ws.Range("E30").Value = "  +2.90%  "

$ws.Range("E31").Value = "  +2.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.76"
$ws.Range("E32").Value = "  +0.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.59"
$ws.Range("E33").Value = "  +0.80%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.44"
$ws.Range("E34").Value = "  +2.94%  "

$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("E36").Value = "  +2.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.07"
$ws.Range("E37").Value = "  +7.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.77"
$ws.Range("E38").Value = "  +3.48%  "

$ws.Range("E39").Value = "  +8.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "126.23"
$ws.Range("E40").Value = "  +4.10%  "

$ws.Range("E41").Value = "  +1.92%  "

$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.24"
$ws.Range("E42").Value = "  +1.53%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.39"
$ws.Range("E43").Value = "  +6.28%  "

$ws.Range("E44").Value = "  +4.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.38"
$ws.Range("E45").Value = "  +8.25%  "

$ws.Range("D46").Value = "2.072.81"
$ws.Range("E46").Value = "  +2.83%  "

$ws.Range("E47").Value = "  +14.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.98"
$ws.Range("E48").Value = "  +8.24%  "

$ws.Range("E49").Value = "  +1.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.42"
$ws.Range("E50").Value = "  +5.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.76"
$ws.Range("E51").Value = "  +1.84%  "

